$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the cell to hold a literal text value (not auto-converted to a
    # number) while keeping the cell's style index unchanged (no "s" attr
    # added to the resulting XML).
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

$ws.Range("D2").Value = "22.409.90"
$ws.Range("E2").Value = "  +0.04%  "

$ws.Range("D3").Value = "1.572.77"
$ws.Range("E3").Value = "  +0.12%  "

$ws.Range("E4").Value = "  +0.23%  "

Set-TextValue "D5" "1.003"
$ws.Range("E5").Value = "  +0.24%  "

Set-TextValue "D6" "290.79"
$ws.Range("E6").Value = "  -0.11%  "

Set-TextValue "D7" "0.3766"
$ws.Range("E7").Value = "  +3.21%  "

Set-TextValue "D8" "49.88"
$ws.Range("E8").Value = "  +0.88%  "

Set-TextValue "D9" "0.3421"
$ws.Range("E9").Value = "  +1.47%  "

Set-TextValue "D10" "1.160"
$ws.Range("E10").Value = "  -0.86%  "

Set-TextValue "D11" "0.07653"
$ws.Range("E11").Value = "  +0.94%  "

$ws.Range("E12").Value = "  +0.23%  "

Set-TextValue "D13" "21.25"
$ws.Range("E13").Value = "  +0.57%  "

Set-TextValue "D14" "6.013"
$ws.Range("E14").Value = "  -0.70%  "

Set-TextValue "D15" "6.934"
$ws.Range("E15").Value = "  +1.03%  "

$ws.Range("D16").Value = "1.572.76"
$ws.Range("E16").Value = "  +0.18%  "

Set-TextValue "D17" "0.00001133"
$ws.Range("E17").Value = "  -0.39%  "

Set-TextValue "D18" "89.92"
$ws.Range("E18").Value = "  +0.97%  "

Set-TextValue "D19" "0.06761"
$ws.Range("E19").Value = "  +0.50%  "

Set-TextValue "D21" "16.81"
$ws.Range("E21").Value = "  +2.19%  "

Set-TextValue "D22" "6.211"
$ws.Range("E22").Value = "  -0.92%  "

Set-TextValue "D23" "12.03"
$ws.Range("E23").Value = "  +0.30%  "

$ws.Range("D24").Value = "22.403.11"

Set-TextValue "D25" "2.420"
$ws.Range("E25").Value = "  +1.14%  "

Set-TextValue "D26" "2.720"
$ws.Range("E26").Value = "  -9.09%  "

Set-TextValue "D27" "20.26"
$ws.Range("E27").Value = "  +1.99%  "

Set-TextValue "D28" "146.64"
$ws.Range("E28").Value = "  +1.64%  "

Set-TextValue "D29" "5.027"
$ws.Range("E29").Value = "  +0.65%  "

Set-TextValue "D30" "126.23"
$ws.Range("E30").Value = "  +0.80%  "

$ws.Range("D31").Value = "1.747.44"
$ws.Range("E31").Value = "  +0.00%  "

Set-TextValue "D32" "6.185"
$ws.Range("E32").Value = "  -1.76%  "

Set-TextValue "D33" "2.010"
$ws.Range("E33").Value = "  +1.90%  "

Set-TextValue "D34" "0.9944"
$ws.Range("E34").Value = "  -5.66%  "

$ws.Range("E35").Value = "  -3.08%  "

Set-TextValue "D36" "0.08604"
$ws.Range("E36").Value = "  +2.05%  "

$ws.Range("E37").Value = "  -0.23%  "

Set-TextValue "D38" "0.2314"
$ws.Range("E38").Value = "  +0.31%  "

Set-TextValue "D39" "0.06573"
$ws.Range("E39").Value = "  +0.73%  "

Set-TextValue "D40" "1.327"
$ws.Range("E40").Value = "  +5.95%  "

$ws.Range("E41").Value = "  -0.84%  "

Set-TextValue "D42" "0.6458"
$ws.Range("E42").Value = "  +1.17%  "

Set-TextValue "D43" "11.53"
$ws.Range("E43").Value = "  -2.56%  "

$ws.Range("E44").Value = "  -2.87%  "

Set-TextValue "D45" "1.002"
$ws.Range("E45").Value = "  +0.27%  "

Set-TextValue "D46" "3.799"
$ws.Range("E46").Value = "  +0.72%  "

$ws.Range("E47").Value = "  -0.20%  "

Set-TextValue "D48" "1.305"
$ws.Range("E48").Value = "  +7.58%  "

Set-TextValue "D49" "2.085"
$ws.Range("E49").Value = "  -2.38%  "

Set-TextValue "D50" "125.40"
$ws.Range("E50").Value = "  +2.43%  "

Set-TextValue "D51" "0.07332"
$ws.Range("E51").Value = "  +0.60%  "
